$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEST_CASES")

# Insert three new columns before the existing column X (TC_CUF_<CODE1>)
# so the sheet gains TC_KIND / TC_SCRIPTING_LANGUAGE / TC_SCRIPT ahead of
# the custom-field / rollup columns that used to start at X.
$ws.Columns("X:Z").Insert()

$ws.Range("X1").Value = "TC_KIND"
$ws.Range("Y1").Value = "TC_SCRIPTING_LANGUAGE"
$ws.Range("Z1").Value = "TC_SCRIPT"

# Match the author's widened columns for the new script-related fields.
$ws.Columns("Y").ColumnWidth = 23.666666666666668
$ws.Columns("Z").ColumnWidth = 10.166666666666666
